$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1: replace the highlighted "password" placeholder text
#    ("___________") with "UML" so students are told to look at the board
#    instead of having the password baked into the deck.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j, 1)
            if ($para.Text -eq "___________") {
                $para.Text = "UML"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Refresh the cached text of every "datetimeFigureOut" style date
#    placeholder (slide master, every slide layout, handout master, notes
#    master) so they show the current date instead of the stale one.
# ---------------------------------------------------------------------------
$longDate = "Friday, March 14, 2025"
$shortDate = "3/14/2025"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $cur = $tr.Text
                if ($cur -match "^\d+/\d+/\d+$") {
                    $tr.Text = $shortDate
                } else {
                    $tr.Text = $longDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# Handout master
if ($p.HasHandoutMaster) {
    Update-DatePlaceholder $p.HandoutMaster.Shapes
}

# Notes master
if ($p.HasNotesMaster) {
    Update-DatePlaceholder $p.NotesMaster.Shapes
}
